# backend/expense_details.xlsx - update expense rows, add a new "Food" row
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Shopping/800 -> Loan/1000 (date moved forward a day)
$ws.Range("A2").Value = "Loan"
$ws.Range("B2").Value = 1000
$ws.Range("C2").Value = 45901.229537037034

# Row 3: Shopping/850 -> Groceries/1000
$ws.Range("A3").Value = "Groceries"
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 45896.229537037034

# Row 4: Rent/1500 -> Shopping/3000
$ws.Range("A4").Value = "Shopping"
$ws.Range("B4").Value = 3000
$ws.Range("C4").Value = 45894.229537037034

# New row 5: Food/1500, picking up the same date style (s="1") used by column C
$ws.Range("A5").Value = "Food"
$ws.Range("B5").Value = 1500
$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C5").Value = 45886.229537037034
